$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "SPRINT 0"

# --- Insert a new row before the existing SUM row (row 14) ---------------
# This pushes the old row 14 (SUM formula) down to row 15, leaving row 14
# empty, exactly like in the target workbook.
$ws.Rows.Item(14).Insert()

# --- Fill the new activity row 13 -----------------------------------------
# Copy formatting (date number format, etc.) from the row above so the new
# cells pick up the same style indices as the rest of the table.
$ws.Range("B12:D12").Copy()
$ws.Range("B13:D13").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B13").Value = 43899
$ws.Range("C13").Value = "Rencontre TB, fin Sprint 0"
$ws.Range("D13").Value = 1

# --- Update the total formula (now living on row 15) ----------------------
$ws.Range("D15").Formula = "=SUM(D3:D13)"

# --- Add the review note on row 17 ----------------------------------------
$ws.Range("C17").Value = "Sprint 0 review, notes de séance en document annexe"

# --- Page setup: portrait / A4-ish paper size for the cahier des charges --
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Restore the active selection shown in the saved file ------------------
$ws.Range("C18").Select() | Out-Null
